$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 70, shifting existing rows 70.. down to 71..
$ws.Rows(70).Insert()

# Populate the newly inserted row 70 with the required values.
$ws.Cells.Item(70, 1).Value = 7
$ws.Cells.Item(70, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(70, 3).Value = "Ñuble"
$ws.Cells.Item(70, 4).Value = 45280
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 100112001
$ws.Cells.Item(70, 7).Value = "Berenjena"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 60
$ws.Cells.Item(70, 11).Value = 11000
$ws.Cells.Item(70, 12).Value = 12000
$ws.Cells.Item(70, 13).Value = 11500
$ws.Cells.Item(70, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 192
$ws.Cells.Item(70, 17).Value = 60
$ws.Cells.Item(70, 18).Value = "Hortaliza"

# Match the date formatting/style used by the other date cells in column D.
$ws.Cells.Item(70, 4).NumberFormat = $ws.Cells.Item(71, 4).NumberFormat
